# Auto-generated edit script applying numeric updates to Kujata_Profits workbook
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H112").Value = 1818.1364
$ws.Range("J112").Value = 1986.6666
$ws.Range("L112").Value = 5959.9998
$ws.Range("N112").Value = -8175.9998
$ws.Range("H121").Value = 897.5
$ws.Range("I121").Value = 0
$ws.Range("J121").Value = 897.5
$ws.Range("K121").Value = 0
$ws.Range("L121").Value = 2692.5
$ws.Range("M121").ClearContents()
$ws.Range("N121").Value = -6186.5
$ws.Range("H125").Value = 2363.2964
$ws.Range("J125").Value = 2440.4666
$ws.Range("L125").Value = 21964.1994
$ws.Range("N125").Value = -26884.1994
$ws.Range("H129").Value = 854
$ws.Range("J129").Value = 869.4400000000001
$ws.Range("L129").Value = 2608.32
$ws.Range("N129").Value = -12608.32
$ws.Range("H132").Value = 7583362
$ws.Range("I132").Value = 8776900
$ws.Range("K132").Value = 26330700
$ws.Range("M132").Value = -26328170
$ws.Range("H135").Value = 41667620
$ws.Range("I135").Value = 564.625
$ws.Range("K135").Value = 5081.625
$ws.Range("M135").Value = -2546.625
$ws.Range("H137").Value = 2618.681
$ws.Range("I137").Value = 2217.4285
$ws.Range("J137").Value = 2942.7693
$ws.Range("K137").Value = 6652.2855
$ws.Range("L137").Value = 8828.3079
$ws.Range("M137").Value = -4102.2855
$ws.Range("N137").Value = -13928.3079
$ws.Range("H138").Value = 2503.4885
$ws.Range("I138").Value = 2315.6667
$ws.Range("J138").Value = 2533.1448
$ws.Range("K138").Value = 6947.000100000001
$ws.Range("L138").Value = 7599.4344
$ws.Range("M138").Value = -1807.000100000001
$ws.Range("N138").Value = -17879.4344
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 10106.527
$ws.Range("I32").Value = 7359.324
$ws.Range("K32").Value = 7359.324
$ws.Range("M32").Value = -7072.324
$ws.Range("H61").Value = 111112860
$ws.Range("I61").Value = 142858540
$ws.Range("K61").Value = 142858540
$ws.Range("M61").Value = -142858328
$ws.Range("H74").Value = 1330.8928
$ws.Range("I74").Value = 938.7826
$ws.Range("J74").Value = 3134.6
$ws.Range("K74").Value = 938.7826
$ws.Range("L74").Value = 3134.6
$ws.Range("M74").Value = -64.7826
$ws.Range("N74").Value = -4882.6
$ws.Range("H77").Value = 1330.8928
$ws.Range("I77").Value = 938.7826
$ws.Range("J77").Value = 3134.6
$ws.Range("K77").Value = 4693.913
$ws.Range("L77").Value = 15673
$ws.Range("M77").Value = -325.9130000000005
$ws.Range("N77").Value = -24409
$ws.Range("H97").Value = 747.4211
$ws.Range("I97").Value = 761.93335
$ws.Range("J97").Value = 693
$ws.Range("K97").Value = 761.93335
$ws.Range("L97").Value = 693
$ws.Range("M97").Value = -265.93335
$ws.Range("N97").Value = -1685
$ws.Range("H110").Value = 203.54546
$ws.Range("I110").Value = 216.5
$ws.Range("K110").Value = 216.5
$ws.Range("M110").Value = 1828.5
$ws.Range("H132").Value = 3706.7585
$ws.Range("I132").Value = 3412.9092
$ws.Range("J132").Value = 4630.2856
$ws.Range("K132").Value = 10238.7276
$ws.Range("L132").Value = 13890.8568
$ws.Range("M132").Value = -7708.7276
$ws.Range("N132").Value = -18950.8568
$ws.Range("H136").Value = 111112860
$ws.Range("I136").Value = 142858540
$ws.Range("K136").Value = 428575620
$ws.Range("M136").Value = -428573070
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H76").Value = 25000
$ws.Range("J76").Value = 25000
$ws.Range("L76").Value = 25000
$ws.Range("N76").Value = -25630
$ws.Range("H79").Value = 25000
$ws.Range("J79").Value = 25000
$ws.Range("L79").Value = 25000
$ws.Range("N79").Value = -27184
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1646.3182
$ws.Range("I31").Value = 1489.8948
$ws.Range("J31").Value = 2637
$ws.Range("K31").Value = 1489.8948
$ws.Range("L31").Value = 2637
$ws.Range("M31").Value = -1194.8948
$ws.Range("N31").Value = -3227
$ws.Range("H34").Value = 1646.3182
$ws.Range("I34").Value = 1489.8948
$ws.Range("J34").Value = 2637
$ws.Range("K34").Value = 1489.8948
$ws.Range("L34").Value = 2637
$ws.Range("M34").Value = -1287.8948
$ws.Range("N34").Value = -3041
$ws.Range("H96").Value = 11450.5
$ws.Range("J96").Value = 11450.5
$ws.Range("L96").Value = 11450.5
$ws.Range("N96").Value = -16942.5
$ws.Range("H132").Value = 2179.5
$ws.Range("I132").Value = 1860.8889
$ws.Range("K132").Value = 5582.6667
$ws.Range("M132").Value = -3052.6667
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H19").Value = 0
$ws.Range("I19").Value = 0
$ws.Range("J19").Value = 0
$ws.Range("K19").Value = 0
$ws.Range("L19").Value = 0
$ws.Range("M19").ClearContents()
$ws.Range("N19").ClearContents()
$ws.Range("H74").Value = 4900
$ws.Range("J74").Value = 4900
$ws.Range("L74").Value = 14700
$ws.Range("N74").Value = -16822
$ws.Range("H77").Value = 4900
$ws.Range("J77").Value = 4900
$ws.Range("L77").Value = 44100
$ws.Range("N77").Value = -54708
$ws.Range("H80").Value = 2872.125
$ws.Range("I80").Value = 995.5
$ws.Range("K80").Value = 2986.5
$ws.Range("M80").Value = -2050.5
$ws.Range("H82").Value = 7740.3125
$ws.Range("J82").Value = 8926.538
$ws.Range("L82").Value = 26779.614
$ws.Range("N82").Value = -27591.614
$ws.Range("H83").Value = 2872.125
$ws.Range("I83").Value = 995.5
$ws.Range("K83").Value = 8959.5
$ws.Range("M83").Value = -4279.5
$ws.Range("H85").Value = 7740.3125
$ws.Range("J85").Value = 8926.538
$ws.Range("L85").Value = 26779.614
$ws.Range("N85").Value = -29587.614
$ws.Range("H92").Value = 542.875
$ws.Range("I92").Value = 544.1818
$ws.Range("J92").Value = 540
$ws.Range("K92").Value = 1632.5454
$ws.Range("L92").Value = 1620
$ws.Range("M92").Value = -384.5454
$ws.Range("N92").Value = -4116
$ws.Range("H112").Value = 66679120
$ws.Range("J112").Value = 76936830
$ws.Range("L112").Value = 230810490
$ws.Range("N112").Value = -230812706
$ws.Range("H131").Value = 25038508
$ws.Range("I131").Value = 333333860
$ws.Range("J131").Value = 41588.49
$ws.Range("K131").Value = 1000001580
$ws.Range("L131").Value = 124765.47
$ws.Range("M131").Value = -999996540
$ws.Range("N131").Value = -134845.47
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H92").Value = 18349
$ws.Range("J92").Value = 18349
$ws.Range("L92").Value = 18349
$ws.Range("N92").Value = -22093
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H62").Value = 39300
$ws.Range("J62").Value = 39300
$ws.Range("L62").Value = 39300
$ws.Range("N62").Value = -40548
$ws.Range("H65").Value = 39300
$ws.Range("J65").Value = 39300
$ws.Range("L65").Value = 117900
$ws.Range("N65").Value = -124140
$ws.Range("H132").Value = 2443.4102
$ws.Range("I132").Value = 2167.9546
$ws.Range("J132").Value = 2799.8823
$ws.Range("K132").Value = 6503.8638
$ws.Range("L132").Value = 8399.6469
$ws.Range("M132").Value = -3973.8638
$ws.Range("N132").Value = -13459.6469
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 7578009.5
$ws.Range("I122").Value = 11907151
$ws.Range("J122").Value = 2011.25
$ws.Range("K122").Value = 35721453
$ws.Range("L122").Value = 6033.75
$ws.Range("M122").Value = -35719003
$ws.Range("N122").Value = -10933.75
$ws.Range("H132").Value = 2222.2888
$ws.Range("I132").Value = 2229.2646
$ws.Range("J132").Value = 2200.7273
$ws.Range("K132").Value = 6687.793799999999
$ws.Range("L132").Value = 6602.1819
$ws.Range("M132").Value = -4157.793799999999
$ws.Range("N132").Value = -11662.1819
$ws.Range("H138").Value = 36435
$ws.Range("J138").Value = 36435
$ws.Range("L138").Value = 36435
$ws.Range("N138").Value = -46715

Write-Output "Applied 201 cell updates"